# "#12 Scemas reshaped on solution slides on presentation"
#
# Reshapes/repositions the 3x2 "schema" picture grid (and the two small
# right-arrow connector shapes between them) on the Result/solution slide.
#
# All target values below were authored in EMU (as stored in the OOXML) and
# converted to points for the PowerPoint object model (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

$emuPerPt = 12700.0

function Set-ShapeRect {
    param(
        $Shape,
        [double]$XEmu,
        [double]$YEmu,
        [double]$CxEmu,
        [double]$CyEmu
    )
    $Shape.Left = $XEmu / $emuPerPt
    $Shape.Top = $YEmu / $emuPerPt
    $Shape.Width = $CxEmu / $emuPerPt
    $Shape.Height = $CyEmu / $emuPerPt
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)

    switch ($shape.Name) {
        "Image 3" {
            Set-ShapeRect -Shape $shape -XEmu 4051883 -YEmu 68130 -CxEmu 4046535 -CyEmu 2846969
        }
        "Image 5" {
            Set-ShapeRect -Shape $shape -XEmu 83890 -YEmu 3103655 -CxEmu 3867652 -CyEmu 1644514
        }
        "Image 7" {
            Set-ShapeRect -Shape $shape -XEmu 4161766 -YEmu 3103655 -CxEmu 3867652 -CyEmu 1644514
        }
        "Image 9" {
            Set-ShapeRect -Shape $shape -XEmu 8230426 -YEmu 3103655 -CxEmu 3867653 -CyEmu 1644514
        }
        "Flèche : droite 4" {
            Set-ShapeRect -Shape $shape -XEmu 3993737 -YEmu 3766659 -CxEmu 151251 -CyEmu 319287
        }
        "Flèche : droite 8" {
            Set-ShapeRect -Shape $shape -XEmu 8073250 -YEmu 3793224 -CxEmu 132009 -CyEmu 319287
        }
    }
}
